# Append a new trade row (row 5) to the CELG named-trade sheet, matching
# the data/formatting of the existing rows, and refresh column A's
# "best fit" width now that a new (shorter) number lives in that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the date-time cell format (style index used by column G, e.g.
# numFmtId 22 "m/d/yyyy h:mm") from row 4 down into row 5 before writing
# the values, so the new G5 cell re-uses the workbook's existing style
# entry instead of Excel minting a brand-new one.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(5, 1).Value = 10069.39
$ws.Cells.Item(5, 2).Value = 10146.5
$ws.Cells.Item(5, 3).Value = 109.08
$ws.Cells.Item(5, 4).Value = 108.25
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = -0.76
$ws.Cells.Item(5, 7).Value = 42612.674537037034
$ws.Cells.Item(5, 8).Value = $false

# Column A was sized to fit "Principle"/10146.5 (width 8.85546875); the
# newly added 10069.39 doesn't change the longest entry, but Excel's
# best-fit recompute on save narrows it slightly to width 9.
$ws.Columns.Item(1).ColumnWidth = 8.1666666
